# Update actions - relative direction
# Row 2 contains per-column attribution values; this updates the
# changed cells to their new values (122 of 189 columns changed).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -0
$ws.Range("B2").Value = -0.1135635860823204
$ws.Range("D2").Value = 0.2627613492971648
$ws.Range("E2").Value = 0.01347081606661644
$ws.Range("F2").Value = -0
$ws.Range("G2").Value = 0
$ws.Range("I2").Value = -0
$ws.Range("J2").Value = -0
$ws.Range("K2").Value = -0.02265144341177558
$ws.Range("L2").Value = -0
$ws.Range("M2").Value = 0.2387385148544099
$ws.Range("N2").Value = 0.005719642659729415
$ws.Range("R2").Value = -0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = -0.1083997172967939
$ws.Range("V2").Value = 0.01236818953665744
$ws.Range("W2").Value = -0.04361436512658663
$ws.Range("AB2").Value = 0
$ws.Range("AC2").Value = -0.07833672262594704
$ws.Range("AD2").Value = 0
$ws.Range("AE2").Value = -0.01171560210610239
$ws.Range("AF2").Value = 0.009298368432744866
$ws.Range("AG2").Value = -0
$ws.Range("AI2").Value = -0
$ws.Range("AJ2").Value = 0
$ws.Range("AK2").Value = 0
$ws.Range("AL2").Value = -0.04486821589741342
$ws.Range("AM2").Value = 0
$ws.Range("AN2").Value = 0.03216264382969419
$ws.Range("AO2").Value = 0.07941566810152632
$ws.Range("AQ2").Value = 0
$ws.Range("AR2").Value = -0
$ws.Range("AT2").Value = 0
$ws.Range("AU2").Value = -0.1587414854446407
$ws.Range("AW2").Value = 0.09496770416775556
$ws.Range("AX2").Value = 0.01841125536205374
$ws.Range("AY2").Value = -0
$ws.Range("BB2").Value = -0
$ws.Range("BC2").Value = -0
$ws.Range("BD2").Value = -0.02837554453630514
$ws.Range("BF2").Value = 0.09823586953220728
$ws.Range("BG2").Value = 0.03785559630282209
$ws.Range("BJ2").Value = -0
$ws.Range("BL2").Value = 0
$ws.Range("BM2").Value = 0.03083793472071189
$ws.Range("BO2").Value = -0.03925015716482984
$ws.Range("BP2").Value = -0.09536144570456224
$ws.Range("BU2").Value = 0
$ws.Range("BV2").Value = -0.03912296556476491
$ws.Range("BX2").Value = 0.02213707306184231
$ws.Range("BY2").Value = -0.01990996932938572
$ws.Range("BZ2").Value = -0
$ws.Range("CB2").Value = 0
$ws.Range("CD2").Value = -0
$ws.Range("CE2").Value = 0.03053216949848069
$ws.Range("CG2").Value = -0.04036857196174517
$ws.Range("CH2").Value = 0.01834765950112722
$ws.Range("CJ2").Value = -0
$ws.Range("CM2").Value = -0
$ws.Range("CN2").Value = -0.01776595600963329
$ws.Range("CP2").Value = 0.03678342491833134
$ws.Range("CQ2").Value = 0.04830583973100779
$ws.Range("CT2").Value = 0
$ws.Range("CU2").Value = -0
$ws.Range("CV2").Value = -0
$ws.Range("CW2").Value = 0.04587107387803037
$ws.Range("CY2").Value = -0.04331952876540518
$ws.Range("CZ2").Value = 0.009187832693228098
$ws.Range("DE2").Value = -0
$ws.Range("DF2").Value = 0.03492830179516355
$ws.Range("DH2").Value = 0.01594690992218723
$ws.Range("DI2").Value = 0.03940935874805038
$ws.Range("DJ2").Value = 0
$ws.Range("DK2").Value = -0
$ws.Range("DL2").Value = -0
$ws.Range("DN2").Value = 0
$ws.Range("DO2").Value = -0.03807029712088379
$ws.Range("DQ2").Value = 0.04801233879067612
$ws.Range("DR2").Value = -0.02779739109799752
$ws.Range("DS2").Value = -0
$ws.Range("DW2").Value = 0
$ws.Range("DX2").Value = -0.07018650623358692
$ws.Range("DY2").Value = -0
$ws.Range("DZ2").Value = -0.01311686632312821
$ws.Range("EA2").Value = -0.03947747523715595
$ws.Range("EB2").Value = 0
$ws.Range("EF2").Value = -0
$ws.Range("EG2").Value = 0.04817555933695059
$ws.Range("EI2").Value = 0.0908645092385791
$ws.Range("EJ2").Value = -0.03158698639538209
$ws.Range("EO2").Value = 0
$ws.Range("EP2").Value = 0.04351795522603605
$ws.Range("EQ2").Value = 0
$ws.Range("ER2").Value = -0.03605048602120684
$ws.Range("ES2").Value = 0.01943489707753068
$ws.Range("ET2").Value = 0
$ws.Range("EU2").Value = -0
$ws.Range("EV2").Value = 0
$ws.Range("EX2").Value = 0
$ws.Range("EY2").Value = 0.04032481743670978
$ws.Range("FA2").Value = -0.0337791203645583
$ws.Range("FB2").Value = 0.01881023979325147
$ws.Range("FD2").Value = -0
$ws.Range("FG2").Value = -0
$ws.Range("FH2").Value = -0.004887850227421824
$ws.Range("FI2").Value = 0
$ws.Range("FJ2").Value = -0.0157778721620328
$ws.Range("FK2").Value = 0.008101020542556516
$ws.Range("FL2").Value = -0
$ws.Range("FN2").Value = -0
$ws.Range("FP2").Value = -0
$ws.Range("FQ2").Value = -0.00516732990964021
$ws.Range("FR2").Value = -0
$ws.Range("FS2").Value = -0.008857138002429203
$ws.Range("FT2").Value = 0.0141670900926394
$ws.Range("FV2").Value = -0
$ws.Range("FW2").Value = -0
$ws.Range("FY2").Value = 0
$ws.Range("FZ2").Value = -0.04489119929027677
$ws.Range("GB2").Value = 0.02588855511003614
$ws.Range("GD2").Value = 0
$ws.Range("GE2").Value = -0
